$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value2 = '65.878.42'
$ws.Range("E2").Value2 = '  +0.95%  '
$ws.Range("D3").Value2 = '3.608.30'
$ws.Range("E3").Value2 = '  +1.99%  '
Set-TextCell "D4" '0.999'
$ws.Range("E4").Value2 = '  -0.21%  '
Set-TextCell "D5" '605.13'
$ws.Range("E5").Value2 = '  -1.00%  '
Set-TextCell "D6" '137.05'
$ws.Range("E6").Value2 = '  -2.26%  '
$ws.Range("D7").Value2 = '3.607.63'
$ws.Range("E7").Value2 = '  +1.98%  '
Set-TextCell "D9" '0.497'
$ws.Range("E9").Value2 = '  +0.98%  '
$ws.Range("E10").Value2 = '  +0.39%  '
$ws.Range("E12").Value2 = '  -0.02%  '
$ws.Range("D13").Value2 = '4.220.40'
$ws.Range("E13").Value2 = '  +1.89%  '
Set-TextCell "D14" '28.10'
$ws.Range("E14").Value2 = '  +3.16%  '
$ws.Range("E15").Value2 = '  -0.01%  '
$ws.Range("D16").Value2 = '3.602.86'
$ws.Range("E16").Value2 = '  +1.33%  '
$ws.Range("E17").Value2 = '  -0.10%  '
$ws.Range("D18").Value2 = '65.021.00'
$ws.Range("E18").Value2 = '  -0.66%  '
Set-TextCell "D19" '10.15'
$ws.Range("E19").Value2 = '  -1.51%  '
Set-TextCell "D20" '14.71'
$ws.Range("E20").Value2 = '  +2.81%  '
Set-TextCell "D21" '5.93'
$ws.Range("E21").Value2 = '  +0.00%  '
Set-TextCell "D22" '399.31'
$ws.Range("E22").Value2 = '  +0.85%  '
Set-TextCell "D23" '0.591'
$ws.Range("E23").Value2 = '  +3.23%  '
$ws.Range("D24").Value2 = '3.751.30'
$ws.Range("E24").Value2 = '  +1.76%  '
Set-TextCell "D25" '74.59'
$ws.Range("E25").Value2 = '  +0.51%  '
$ws.Range("E26").Value2 = '  -0.01%  '
Set-TextCell "D27" '0.0000119'
$ws.Range("E27").Value2 = '  +1.93%  '
Set-TextCell "D28" '8.19'
$ws.Range("E28").Value2 = '  +4.26%  '
$ws.Range("E29").Value2 = '  +30.09%  '
Set-TextCell "D30" '2.41'
$ws.Range("E30").Value2 = '  +4.58%  '
Set-TextCell "D31" '8.70'
$ws.Range("E31").Value2 = '  +4.55%  '
Set-TextCell "D32" '0.999'
$ws.Range("E32").Value2 = '  -0.17%  '
$ws.Range("D33").Value2 = '3.605.49'
$ws.Range("E33").Value2 = '  +1.36%  '
$ws.Range("E34").Value2 = '  +3.78%  '
$ws.Range("E35").Value2 = '  +0.95%  '
$ws.Range("E36").Value2 = '  -0.02%  '
Set-TextCell "D37" '5.41'
$ws.Range("E37").Value2 = '  +7.82%  '
$ws.Range("B38").Value2 = 'ImmutableX'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D38" '1.61'
$ws.Range("E38").Value2 = '  +3.08%  '
$ws.Range("B39").Value2 = 'Aptos'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D39" '7.13'
$ws.Range("E39").Value2 = '  +1.49%  '
Set-TextCell "D40" '171.99'
$ws.Range("E40").Value2 = '  +1.23%  '
$ws.Range("E41").Value2 = '  +2.63%  '
Set-TextCell "D42" '0.844'
$ws.Range("E42").Value2 = '  +1.58%  '
Set-TextCell "D43" '26.18'
$ws.Range("E43").Value2 = '  -0.65%  '
Set-TextCell "D44" '43.46'
$ws.Range("E44").Value2 = '  +1.19%  '
$ws.Range("E45").Value2 = '  +3.83%  '
Set-TextCell "D46" '4.56'
$ws.Range("E46").Value2 = '  +2.59%  '
$ws.Range("E47").Value2 = '  -0.28%  '
$ws.Range("E48").Value2 = '  +1.31%  '
Set-TextCell "D49" '7.10'
$ws.Range("E49").Value2 = '  +4.03%  '
$ws.Range("D50").Value2 = '2.470.22'
$ws.Range("E50").Value2 = '  +0.50%  '
$ws.Range("E51").Value2 = '  +3.35%  '
